# "Hoan thien Ngoai Tru" - update reception (Tiep nhan) test data on the
# "Data" sheet: new patient record id and insurance card number, plus the
# resulting column widening for the columns that now hold wider values, and
# move the active selection to where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2 test data updates
$ws.Range("A2").Value = 3012
$ws.Range("E2").Value = 46200608012

# Columns E (InsCardNo) and F got wider thanks to the longer values they now
# hold - widen them to fit (approx. 12 and 10 characters respectively).
$ws.Columns.Item(5).ColumnWidth = 11.17
$ws.Columns.Item(6).ColumnWidth = 9.17

# Leave the selection where the editor ended up working
$ws.Range("M14").Select()
